$wb = $excel.ActiveWorkbook

# Sheet 2: "Tipo de questoes aceitas" gets populated with a new row of
# question data (mirrors the existing "discursiva" row but for objetiva03)
# and becomes the active / selected sheet.
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A6").Value = "objetiva03"
$ws2.Range("B6").Value = 2.5
$ws2.Range("C6").Value = "Considere as afirmações abaixo e assinale  a alternativa que indique as afirmativa(s) INCORRETA(S):"
$ws2.Range("D6").Value = "AAAA"
$ws2.Range("E6").Value = "BBBB"
$ws2.Range("F6").Value = "CCCCC"
$ws2.Range("G6").Value = "R: I e II"
$ws2.Range("H6").Value = "R: II e III"
$ws2.Range("I6").Value = "R: I e III"
$ws2.Range("J6").Value = "R: Nenhuma"
$ws2.Range("K6").Value = "R: Todas"

# G6 picked up the wrong style from auto-fill; re-apply the formatting used
# by the rest of the row (same as G5) without touching its new value.
$ws2.Range("G5").Copy()
$ws2.Range("G6").PasteSpecial(-4122)

# Row 6 grows to fit the new wrapped text, row 7 reverts to an
# auto-calculated (non-custom) height.
$ws2.Rows.Item(6).RowHeight = 67
$ws2.Rows.Item(7).AutoFit()

# Make sheet 2 the active sheet/tab and select the newly filled cell.
$ws2.Activate()
$ws2.Range("A6").Select()
